$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.045.06'
$ws.Range('E2').Value = '  +3.19%  '
$ws.Range('D3').Value = '1.688.45'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.54'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.519'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.10'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E9').Value = '  +1.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0626'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.33%  '
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').Value = '1.925.76'
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').Value = '1.689.58'
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.88'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '250.62'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.10%  '
$ws.Range('D18').Value = '27.994.03'
$ws.Range('E18').Value = '  +3.05%  '
$ws.Range('D19').Value = '0.0₃0743'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.96%  '
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.55'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.05'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.94%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.60'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.35'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.49'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.98%  '
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  +6.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0504'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('E33').Value = '  -2.02%  '
$ws.Range('D34').Value = '1.426.06'
$ws.Range('E34').Value = '  -7.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.62'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.944'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.592'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('E40').Value = '  -3.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.43'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.50'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.50%  '
$ws.Range('D44').Value = '1.833.23'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.24'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.92%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.796'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.55%  '
$ws.Range('E47').Value = '  +5.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '89.41'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.78%  '
$ws.Range('D49').Value = '0.0₆0111'
$ws.Range('E49').Value = '  -0.93%  '
$ws.Range('E50').Value = '  -0.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.85'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.45%  '
